# Reprocess the metadata sheet with the newly curated dimensions:
# - "lugares-de-importancia-comunitaria", "espacios-naturales-protegidos" and
#   "zonas-de-especial-proteccion-para-las-aves" stop being curated dimensions
#   (iaest-dimension:* / dim / skos:Concept) and become plain measures
#   (iaest-measure:* / medida / xsd:int) -> columns G, I, L
# - "municipio-nombre" becomes the curated geo-area dimension instead
#   (sdmx-dimension:refArea / dim / URI-Municipio) -> column H
# - the mapping workbooks that used to back the now-dropped dimensions are
#   no longer referenced -> row 5 cells G5, I5, L5

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G: lugares-de-importancia-comunitaria (was a curated dimension)
$ws.Range("G2").Value = "iaest-measure:lugares-de-importancia-comunitaria"
$ws.Range("G3").Value = "medida"
$ws.Range("G4").Value = "xsd:int"
$ws.Range("G5").Value = ""

# Column H: municipio-nombre (becomes the curated geo-area dimension)
$ws.Range("H2").Value = "sdmx-dimension:refArea"
$ws.Range("H3").Value = "dim"
$ws.Range("H4").Value = "URI-Municipio"

# Column I: espacios-naturales-protegidos (was a curated dimension)
$ws.Range("I2").Value = "iaest-measure:espacios-naturales-protegidos"
$ws.Range("I3").Value = "medida"
$ws.Range("I4").Value = "xsd:int"
$ws.Range("I5").Value = ""

# Column L: zonas-de-especial-proteccion-para-las-aves (was a curated dimension)
$ws.Range("L2").Value = "iaest-measure:zonas-de-especial-proteccion-para-las-aves"
$ws.Range("L3").Value = "medida"
$ws.Range("L4").Value = "xsd:int"
$ws.Range("L5").Value = ""
